$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Hallitusbotti on kovis! @ 11.4.2017, 12:57"
$ws.Range("A4").Value = "TIEDÄN HALLITUSBOTIN SALASANAN ENKÄ OLE HALLITUKSESSA! Paska hallitus.. @ 11.4.2017, 13:25"
